$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Coal): Capacity_MW 100 -> 400, Fuel_Cost_USDperMMBtu 2.41 -> 2.5499999999999998 ---
$ws.Range("B2").Value = 400
$ws.Range("D2").Value = 2.55

# --- Column I (Carbon_tonsperMWH): change divisor from 2000 to 2204.62, as a shared formula I3:I6 ---
$ws.Range("I2").Formula = "=H2*C2/2204.62"
$ws.Range("I3:I6").Formula = "=H3*C3/2204.62"

# The formula fill into I6 also carried the formatting used by the rest of the
# formula column (I2:I5), overriding I6's previous one-off border style.
$ws.Range("I5").Copy() | Out-Null
$ws.Range("I6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Row 3 (Natural Gas): Fuel_Cost_USDperMMBtu 6.45 -> 3.03 ---
$ws.Range("D3").Value = 3.03

# --- Row 4 (Nuclear): Capacity_MW 100 -> 1000, Fuel_Cost_USDperMMBtu blank -> 0.97, Fuel_Cost_USDperMWH becomes a formula ---
$ws.Range("B4").Value = 1000
$ws.Range("D4").Value = 0.97
$ws.Range("E4").Formula = "=D4*C4"

# --- Row 5 (Wind): Heat_Rate_MMBtu_perMWH blank -> 0 ---
$ws.Range("C5").Value = 0

# --- Row 6 (Solar): Heat_Rate_MMBtu_perMWH blank -> 0 ---
$ws.Range("C6").Value = 0

# --- Update the active selection to match the author's final cursor position ---
$ws.Range("I5").Select() | Out-Null

$wb.Save()
